{"js": "// Reworks the \"Regras b\u00e1sicas\" bullet list:\n//  - item 1 (\"Os requisitos devem ser escritos de forma clara e objetiva.\")\n//    is replaced by the old item 2's text, split into two runs on the colon\n//  - item 2 (\"Definir apenas um requisito de cada vez.\") becomes old item 3's\n//    text, split into two runs\n//  - item 3 (\"Evitar requisitos muito extensos.\") becomes old item 4's text,\n//    split into two runs\n//  - item 4 (\"Para o uso de termos t\u00e9cnicos deve existir um gloss\u00e1rio.\")\n//    becomes the new \"Evite palavras ou frases...\" text, split into two runs\n//  - item 5 (\"N\u00e3o utilizar palavras que causem uma ambiguidade ao\n//    requisito.\") is removed entirely\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the five bullet paragraphs by their (old) text so the script does\n// not depend on hard-coded indices.\nconst wanted = [\n  \"Os requisitos devem ser escritos de forma clara e objetiva.\",\n  \"Definir apenas um requisito de cada vez.\",\n  \"Evitar requisitos muito extensos.\",\n  \"Para o uso de termos t\u00e9cnicos deve existir um gloss\u00e1rio.\",\n  \"N\u00e3o utilizar palavras que causem uma ambiguidade ao requisito.\",\n];\n\nconst found = {};\nfor (const item of paragraphs.items) {\n  const t = (item.text || \"\").trim();\n  for (const w of wanted) {\n    if (t === w) {\n      found[w] = item;\n    }\n  }\n}\n\n// New content for each of the first four bullets: [firstRunText, secondRunText]\nconst replacements = [\n  {\n    key: wanted[0],\n    run1: \"Definir apenas um requisito de cada vez\",\n    run2: \": ou seja, o requisito deve ser respons\u00e1vel por apenas uma fun\u00e7\u00e3o, caso seja necess\u00e1rio mais de uma, separar em dois ou mais requisitos;\",\n  },\n  {\n    key: wanted[1],\n    run1: \"Evitar requisitos muito extensos\",\n    run2: \": evitar frases grandes com palavras que n\u00e3o acrescentem nada para o entendimento do requisito;\",\n  },\n  {\n    key: wanted[2],\n    run1: \"Para o uso de termos t\u00e9cnicos deve existir um gloss\u00e1rio\",\n    run2: \";\",\n  },\n  {\n    key: wanted[3],\n    run1: \"Evite palavras ou frases como: e, ou, somente se, exceto, se necess\u00e1rio, mas, contudo, entretanto, usualmente, geralmente, frequentemente, tipicamente, amig\u00e1vel, vers\u00e1til, flex\u00edvel, aproximadamente, t\u00e3o logo quanto poss\u00edvel, talvez, \",\n    run2: \"provavelmente etc.\",\n  },\n];\n\nfor (const r of replacements) {\n  const para = found[r.key];\n  if (!para) {\n    throw new Error(\"Could not locate paragraph with text: \" + r.key);\n  }\n  const range = para.getRange();\n  // Replace the paragraph's whole text with the first run's text, keeping\n  // the original run formatting (sz/szCs) intact.\n  range.insertText(r.run1, Word.InsertLocation.replace);\n  await context.sync();\n\n  // Insert the second run right after the first one. Nudging the\n  // formatting (set true, then back to false, each with its own sync)\n  // keeps the two runs from being re-merged into a single <w:r> during\n  // serialization, matching the two separate <w:r> elements produced by\n  // the original edit.\n  const run2Range = range.insertText(r.run2, Word.InsertLocation.after);\n  run2Range.font.bold = true;\n  await context.sync();\n  run2Range.font.bold = false;\n  await context.sync();\n}\n\n// Remove the fifth bullet entirely.\nconst lastPara = found[wanted[4]];\nif (!lastPara) {\n  throw new Error(\"Could not locate paragraph with text: \" + wanted[4]);\n}\nlastPara.delete();\nawait context.sync();\n", "ps1": "# Reworks the \"Regras b\u00e1sicas\" bullet list:\n#  - item 1 (\"Os requisitos devem ser escritos de forma clara e objetiva.\")\n#    is replaced by the old item 2's text, split into two runs on the colon\n#  - item 2 (\"Definir apenas um requisito de cada vez.\") becomes old item 3's\n#    text, split into two runs\n#  - item 3 (\"Evitar requisitos muito extensos.\") becomes old item 4's text,\n#    split into two runs\n#  - item 4 (\"Para o uso de termos t\u00e9cnicos deve existir um gloss\u00e1rio.\")\n#    becomes the new \"Evite palavras ou frases...\" text, split into two runs\n#  - item 5 (\"N\u00e3o utilizar palavras que causem uma ambiguidade ao\n#    requisito.\") is removed entirely\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphText($para) {\n    $raw = $para.Range.Text\n    return $raw.TrimEnd([char]13, [char]10)\n}\n\n# Locate the five target bullet paragraphs by their current text so the\n# script does not depend on hard-coded indices.\n$targets = @(\n    \"Os requisitos devem ser escritos de forma clara e objetiva.\",\n    \"Definir apenas um requisito de cada vez.\",\n    \"Evitar requisitos muito extensos.\",\n    \"Para o uso de termos t\u00e9cnicos deve existir um gloss\u00e1rio.\",\n    \"N\u00e3o utilizar palavras que causem uma ambiguidade ao requisito.\"\n)\n\n$foundParas = @{}\nforeach ($p in $d.Paragraphs) {\n    $t = Get-ParagraphText $p\n    foreach ($target in $targets) {\n        if ($t -eq $target) {\n            $foundParas[$target] = $p\n        }\n    }\n}\n\n# New content for the first four bullets: run1 (keeps old formatting) then\n# run2 (inserted right after, same formatting, kept as its own <w:r>).\n$replacements = @(\n    @{ Key = $targets[0]; Run1 = \"Definir apenas um requisito de cada vez\"; Run2 = \": ou seja, o requisito deve ser respons\u00e1vel por apenas uma fun\u00e7\u00e3o, caso seja necess\u00e1rio mais de uma, separar em dois ou mais requisitos;\" },\n    @{ Key = $targets[1]; Run1 = \"Evitar requisitos muito extensos\"; Run2 = \": evitar frases grandes com palavras que n\u00e3o acrescentem nada para o entendimento do requisito;\" },\n    @{ Key = $targets[2]; Run1 = \"Para o uso de termos t\u00e9cnicos deve existir um gloss\u00e1rio\"; Run2 = \";\" },\n    @{ Key = $targets[3]; Run1 = \"Evite palavras ou frases como: e, ou, somente se, exceto, se necess\u00e1rio, mas, contudo, entretanto, usualmente, geralmente, frequentemente, tipicamente, amig\u00e1vel, vers\u00e1til, flex\u00edvel, aproximadamente, t\u00e3o logo quanto poss\u00edvel, talvez, \"; Run2 = \"provavelmente etc.\" }\n)\n\nforeach ($r in $replacements) {\n    $para = $foundParas[$r.Key]\n    if ($null -eq $para) {\n        throw (\"Could not locate paragraph with text: \" + $r.Key)\n    }\n    $range = $para.Range\n\n    # Replace the whole paragraph text (minus its trailing mark) with the\n    # first run's text; Range.Text assignment keeps the existing run\n    # formatting (sz/szCs) intact.\n    $range.Text = $r.Run1\n\n    # Insert the second run right before the paragraph mark so it stays\n    # inside this paragraph.\n    $insertionPoint = $d.Range($range.End - 1, $range.End - 1)\n    $insertionPoint.InsertAfter($r.Run2)\n\n    # Nudge formatting (set then restore) on just the newly inserted text so\n    # it keeps its own <w:r> element instead of being re-merged with the\n    # first run during serialization - this mirrors the two separate <w:r>\n    # elements produced by the original edit.\n    $insertionPoint.Bold = 1\n    $insertionPoint.Bold = 0\n}\n\n# Remove the fifth bullet entirely (paragraph + its paragraph mark).\n$lastPara = $foundParas[$targets[4]]\nif ($null -eq $lastPara) {\n    throw (\"Could not locate paragraph with text: \" + $targets[4])\n}\n$lastPara.Range.Delete()\n"}
